$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three new rows (270-272) repeat the same metric values as row 269,
# only the date serial in column A advances by one day each row.
$values = @(116.4121952, 0.00170247, 0.008850780000000001, 0.06933635, 12792.90181321, 465.80531254, 0.24, 1.7904431, 485.38834923)

$lastExistingRow = 269
$startRow = 270
$startDate = 45826
$numNewRows = 3

for ($i = 0; $i -lt $numNewRows; $i++) {
    $row = $startRow + $i
    $date = $startDate + $i

    # Copy the formatting (number format, alignment, style) from column A of
    # the last existing row so the new date cell reuses the same style index
    # instead of minting a brand-new one.
    $ws.Range("A$lastExistingRow").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $date

    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 2).Value = $values[$c]
    }
}

$excel.CutCopyMode = $false
